$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.745.25'
$ws.Range("E2").Value = '  -4.58%  '
$ws.Range("D3").Value = '1.718.43'
$ws.Range("E3").Value = '  -5.45%  '
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").Value = '223.81'
$ws.Range("E5").Value = '  -4.03%  '
$ws.Range("D6").Value = '0.5600'
$ws.Range("E6").Value = '  -5.01%  '
$ws.Range("D7").Value = '1.004'
$ws.Range("E7").Value = '  -0.25%  '
$ws.Range("D8").Value = '0.2685'
$ws.Range("E8").Value = '  -1.54%  '
$ws.Range("D9").Value = '22.53'
$ws.Range("E9").Value = '  -1.53%  '
$ws.Range("D10").Value = '0.06496'
$ws.Range("E10").Value = '  -4.35%  '
$ws.Range("D11").Value = '0.07490'
$ws.Range("E11").Value = '  -0.43%  '
$ws.Range("D12").Value = '1.723.17'
$ws.Range("E12").Value = '  -4.76%  '
$ws.Range("D13").Value = '4.626'
$ws.Range("E13").Value = '  -0.17%  '
$ws.Range("D14").Value = '0.5887'
$ws.Range("E14").Value = '  -5.26%  '
$ws.Range("D15").Value = '1.953.66'
$ws.Range("E15").Value = '  -5.30%  '
$ws.Range("D16").Value = '72.98'
$ws.Range("E16").Value = '  -2.29%  '
$ws.Range("D17").Value = '0.000008481'
$ws.Range("E17").Value = '  -10.14%  '
$ws.Range("D18").Value = '27.720.14'
$ws.Range("E18").Value = '  -3.95%  '
$ws.Range("D19").Value = '5.214'
$ws.Range("E19").Value = '  -3.72%  '
$ws.Range("E20").Value = '  -0.24%  '
$ws.Range("D21").Value = '11.11'
$ws.Range("E21").Value = '  -2.30%  '
$ws.Range("D22").Value = '198.12'
$ws.Range("E22").Value = '  -4.64%  '
$ws.Range("D23").Value = '6.470'
$ws.Range("E23").Value = '  -4.19%  '
$ws.Range("E24").Value = '  -0.21%  '
$ws.Range("D25").Value = '149.02'
$ws.Range("E25").Value = '  -3.44%  '
$ws.Range("D26").Value = '7.877'
$ws.Range("E26").Value = '  +1.24%  '
$ws.Range("D27").Value = '0.1204'
$ws.Range("E27").Value = '  -4.65%  '
$ws.Range("D28").Value = '15.89'
$ws.Range("E28").Value = '  -2.11%  '
$ws.Range("D29").Value = '1.362'
$ws.Range("E29").Value = '  -3.38%  '
$ws.Range("D30").Value = '0.06047'
$ws.Range("E30").Value = '  -4.59%  '
$ws.Range("D31").Value = '1.378'
$ws.Range("E31").Value = '  -3.55%  '
$ws.Range("D32").Value = '3.637'
$ws.Range("E32").Value = '  -1.65%  '
$ws.Range("D33").Value = '3.642'
$ws.Range("E33").Value = '  -0.93%  '
$ws.Range("D34").Value = '1.658'
$ws.Range("E34").Value = '  -1.67%  '
$ws.Range("D35").Value = '1.020'
$ws.Range("E35").Value = '  -2.69%  '
$ws.Range("D36").Value = '0.6357'
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("D37").Value = '2.422'
$ws.Range("E37").Value = '  -4.40%  '
$ws.Range("D38").Value = '2.678'
$ws.Range("E38").Value = '  -2.41%  '
$ws.Range("D39").Value = '0.01652'
$ws.Range("E39").Value = '  -3.07%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '6.097'
$ws.Range("E40").Value = '  -4.49%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '1.104.21'
$ws.Range("E41").Value = '  -2.51%  '
$ws.Range("D42").Value = '0.8720'
$ws.Range("E42").Value = '  +1.10%  '
$ws.Range("E43").Value = '  -0.18%  '
$ws.Range("D44").Value = '98.75'
$ws.Range("E44").Value = '  -1.33%  '
$ws.Range("D45").Value = '1.867.22'
$ws.Range("E45").Value = '  -5.29%  '
$ws.Range("D46").Value = '58.38'
$ws.Range("E46").Value = '  -2.89%  '
$ws.Range("D47").Value = '0.00000000109'
$ws.Range("E47").Value = '  -3.62%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.05354'
$ws.Range("E48").Value = '  -2.31%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '1.524'
$ws.Range("E49").Value = '  -2.89%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '8.131'
$ws.Range("E50").Value = '  -0.90%  '
$ws.Range("D51").Value = '0.4400'
$ws.Range("E51").Value = '  -3.01%  '
